$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting so that
# numeric-looking values (e.g. "0.998", "225.65") are not auto-converted to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.165.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.021.45"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.65"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.608"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.98"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.378"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0786"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -6.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.314.82"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.23"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.21"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.740"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.18"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.027.55"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.076.77"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.85"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0815"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "222.64"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.73%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.44"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.17"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -6.77%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.33"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.21"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.128"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.69"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.32"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.117"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.51"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0612"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.46"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.35"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.38%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.56"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.62%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.10"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.470.83"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0215"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.12%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "95.51"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.55"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0909"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.75"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.13"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.22"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.204.89"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.58"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -9.74%  "
